# Fruta / hortaliza, semanal
# Insert 3 new daily price records for "Femacal de La Calera" (Arándano (blue))
# at rows 63-65, pushing the previously-existing rows 63..126 down to 66..129.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 63 (formats/styles of the
# row above are carried down automatically by Excel, matching the existing
# "D" column date-style column).
$ws.Range("A63:A65").EntireRow.Insert()

# --- New row 63 ---
$ws.Range("A63").Value = 3
$ws.Range("B63").Value = "Femacal de La Calera"
$ws.Range("C63").Value = "Coquimbo"
$ws.Range("D63").Value = 44494
$ws.Range("E63").Value = 5
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100101
$ws.Range("H63").Value = "Berries"
$ws.Range("I63").Value = 100101001
$ws.Range("J63").Value = "Arándano (blue)"
$ws.Range("K63").Value = "Sin especificar"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 110
$ws.Range("N63").Value = 10000
$ws.Range("O63").Value = 11000
$ws.Range("P63").Value = 10545
$ws.Range("Q63").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R63").Value = "Provincia de Quillota"
$ws.Range("S63").Value = 7030
$ws.Range("T63").Value = 1.5

# --- New row 64 ---
$ws.Range("A64").Value = 3
$ws.Range("B64").Value = "Femacal de La Calera"
$ws.Range("C64").Value = "Coquimbo"
$ws.Range("D64").Value = 44494
$ws.Range("E64").Value = 5
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100101
$ws.Range("H64").Value = "Berries"
$ws.Range("I64").Value = 100101001
$ws.Range("J64").Value = "Arándano (blue)"
$ws.Range("K64").Value = "Sin especificar"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 60
$ws.Range("N64").Value = 10000
$ws.Range("O64").Value = 10000
$ws.Range("P64").Value = 10000
$ws.Range("Q64").Value = "$/bandeja 2 kilos"
$ws.Range("R64").Value = "Provincia de Quillota"
$ws.Range("S64").Value = 5000
$ws.Range("T64").Value = 2

# --- New row 65 ---
$ws.Range("A65").Value = 3
$ws.Range("B65").Value = "Femacal de La Calera"
$ws.Range("C65").Value = "Coquimbo"
$ws.Range("D65").Value = 44494
$ws.Range("E65").Value = 5
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100101
$ws.Range("H65").Value = "Berries"
$ws.Range("I65").Value = 100101001
$ws.Range("J65").Value = "Arándano (blue)"
$ws.Range("K65").Value = "Sin especificar"
$ws.Range("L65").Value = "Segunda"
$ws.Range("M65").Value = 54
$ws.Range("N65").Value = 8000
$ws.Range("O65").Value = 8000
$ws.Range("P65").Value = 8000
$ws.Range("Q65").Value = "$/bandeja 2 kilos"
$ws.Range("R65").Value = "Provincia de Quillota"
$ws.Range("S65").Value = 4000
$ws.Range("T65").Value = 2
